# 06组项目计划表 - add "第九周四" status block (rows 99-108), fill in
# completion percentages for the "第九周一" block (rows 91-95), and
# move the active view down to the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Fill in the completion percentages (column C) for the existing
#    "日期：2018.10.15 第九周一" block (rows 91-95).
# ---------------------------------------------------------------------
$ws.Range("C91").Value = 0.6
$ws.Range("C92").Value = 0.54
$ws.Range("C93").Value = 0.94
$ws.Range("C94").Value = 0.85
$ws.Range("C95").Value = 0.85

# ---------------------------------------------------------------------
# 2. Duplicate the formatting of the previous week block (rows 89-98)
#    down onto the new block (rows 99-108), then overwrite the cell
#    contents for the new "第九周四" entries.
# ---------------------------------------------------------------------
$ws.Range("A89:D98").Copy() | Out-Null
$ws.Range("A99").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Section header / date row
$ws.Range("A99").Value = "日期：2018.10.15 第九周四"
$ws.Range("A99:D99").Merge()

# Column header row
$ws.Range("A100").Value = "组员"
$ws.Range("B100").Value = "计划内容"
$ws.Range("C100").Value = "完成情况"
$ws.Range("D100").Value = "备注"

# Member rows
$ws.Range("A101").Value = "李福森"
$ws.Range("B101").Value = "内容:编写web端普通用户的ui代码"

$ws.Range("A102").Value = "陈添楠"
$ws.Range("B102").Value = "内容:编写web端普通用户的js代码"

$ws.Range("A103").Value = "黄龙强"
$ws.Range("B103").Value = "内容:编写android ui代码"

$ws.Range("A104").Value = "邱培松"
$ws.Range("B104").Value = "内容:编写网络交互处理模块"

$ws.Range("A105").Value = "王一鸣"
$ws.Range("B105").Value = "内容:找android和web端界面ui素材"

# Row 106 stays blank (spacer row), matches the copied template.

# Summary row
$ws.Range("A107").Value = "总结："
$ws.Range("A107:D108").Merge()

# Row 108 stays blank (bottom border row of the summary box).

# ---------------------------------------------------------------------
# 3. Scroll / select so the newly added block is in view, matching the
#    author's last on-screen position.
# ---------------------------------------------------------------------
$ws.Range("C101").Select()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
